# Fix workbook headers + add PO Forecast sheet

$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item(1)
$wsMonthly = $wb.Worksheets.Item(2)

# 1) Rename the "Requested quantity" headers on the existing sheets.
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# 2) Add a new "PO Forecast" sheet after "Monthly Trend".
$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Reuse the existing header / date formatting (style ids) from the
# "Weekly Quantity" sheet instead of creating brand-new style entries.
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A11").PasteSpecial(-4122)

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Data rows
$wsForecast.Range("A2").Value = 44934.99999999999
$wsForecast.Range("B2").Value = 4
$wsForecast.Range("C2").Value = 4.000015528454801
$wsForecast.Range("D2").Value = 4.000015528884949
$wsForecast.Range("A3").Value = 44941.99999999999
$wsForecast.Range("B3").Value = 20
$wsForecast.Range("C3").Value = 20.00001552855725
$wsForecast.Range("D3").Value = 20.00001552900017
$wsForecast.Range("A4").Value = 44948.99999999999
$wsForecast.Range("B4").Value = 36
$wsForecast.Range("C4").Value = 36.00001537120212
$wsForecast.Range("D4").Value = 36.00001568727066
$wsForecast.Range("A5").Value = 44955.99999999999
$wsForecast.Range("B5").Value = 52
$wsForecast.Range("C5").Value = 52.00001500044697
$wsForecast.Range("D5").Value = 52.00001600427091
$wsForecast.Range("A6").Value = 44962.99999999999
$wsForecast.Range("B6").Value = 68
$wsForecast.Range("C6").Value = 68.00001451900633
$wsForecast.Range("D6").Value = 68.00001644240756
$wsForecast.Range("A7").Value = 44969.99999999999
$wsForecast.Range("B7").Value = 84
$wsForecast.Range("C7").Value = 84.00001392986637
$wsForecast.Range("D7").Value = 84.00001697083087
$wsForecast.Range("A8").Value = 44976.99999999999
$wsForecast.Range("B8").Value = 100
$wsForecast.Range("C8").Value = 100.0000133342487
$wsForecast.Range("D8").Value = 100.0000175611534
$wsForecast.Range("A9").Value = 44983.99999999999
$wsForecast.Range("B9").Value = 116
$wsForecast.Range("C9").Value = 116.0000125918984
$wsForecast.Range("D9").Value = 116.0000182271143
$wsForecast.Range("A10").Value = 44990.99999999999
$wsForecast.Range("B10").Value = 132
$wsForecast.Range("C10").Value = 132.0000118478337
$wsForecast.Range("D10").Value = 132.0000190061631
$wsForecast.Range("A11").Value = 44997.99999999999
$wsForecast.Range("B11").Value = 148
$wsForecast.Range("C11").Value = 148.0000110901495
$wsForecast.Range("D11").Value = 148.000019749539

# Leave the first sheet selected/active, matching the original workbook.
$wsWeekly.Select()
